$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("IK106")
$ws.Range("E2").Value = 45096.49076744096
$ws.Range("E3").Value = 45096.49076744096
$ws.Range("E4").Value = 45096.49076744096
$ws.Range("E5").Value = 45096.49076744096

$ws = $wb.Worksheets.Item("Q100")
$ws.Range("E2").Value = 45096.48779674769
$ws.Range("E3").Value = 45096.48779690972
$ws.Range("E4").Value = 45096.48779674769
$ws.Range("E5").Value = 45096.48779674769

$ws = $wb.Worksheets.Item("Q200")
$ws.Range("E2").Value = 45096.48796245371
$ws.Range("E3").Value = 45096.48796245371
$ws.Range("E4").Value = 45096.48796245371

$ws = $wb.Worksheets.Item("Q300")
$ws.Range("E2").Value = 45096.48813849537
$ws.Range("E3").Value = 45096.48813861111
$ws.Range("E4").Value = 45096.48813861111
$ws.Range("E5").Value = 45096.48813849537

$ws = $wb.Worksheets.Item("Q400")
$ws.Range("E2").Value = 45096.48839530093
$ws.Range("E3").Value = 45096.48839530093
$ws.Range("E4").Value = 45096.48839552084
$ws.Range("E5").Value = 45096.48839552084

$ws = $wb.Worksheets.Item("Q500")
$ws.Range("E2").Value = 45096.48860040509
$ws.Range("E3").Value = 45096.48860040509
$ws.Range("E4").Value = 45096.4886006713
$ws.Range("E5").Value = 45096.48860040509
$ws.Range("E6").Value = 45096.48860040509

$ws = $wb.Worksheets.Item("Q600")
$ws.Range("E2").Value = 45096.48873192129
$ws.Range("E3").Value = 45096.48873192129
$ws.Range("E4").Value = 45096.48873192129
$ws.Range("E5").Value = 45096.48873192129

$ws = $wb.Worksheets.Item("P100")
$ws.Range("E2").Value = 45096.48891115741
$ws.Range("E3").Value = 45096.48891115741
$ws.Range("E4").Value = 45096.48891115741
$ws.Range("E5").Value = 45096.48891115741

$ws = $wb.Worksheets.Item("P200")
$ws.Range("E2").Value = 45096.48905635416
$ws.Range("E3").Value = 45096.48905635416
$ws.Range("E4").Value = 45096.4890566088
$ws.Range("E5").Value = 45096.48905635416
$ws.Range("E6").Value = 45096.48905635416

$ws = $wb.Worksheets.Item("P300")
$ws.Range("E2").Value = 45096.48926726852
$ws.Range("E3").Value = 45096.48926726852
$ws.Range("E4").Value = 45096.48926726852
$ws.Range("E5").Value = 45096.48926726852
$ws.Range("E6").Value = 45096.48926726852
$ws.Range("E7").Value = 45096.48926726852
$ws.Range("E8").Value = 45096.48926726852

$ws = $wb.Worksheets.Item("P400")
$ws.Range("E2").Value = 45096.48950741898
$ws.Range("E3").Value = 45096.48950741898
$ws.Range("E4").Value = 45096.48950741898
$ws.Range("E5").Value = 45096.4895075926
$ws.Range("E6").Value = 45096.48950741898

$ws = $wb.Worksheets.Item("P500")
$ws.Range("E2").Value = 45096.48980850694
$ws.Range("E3").Value = 45096.48980850694
$ws.Range("E4").Value = 45096.48980850694
$ws.Range("E5").Value = 45096.48980870371

$ws = $wb.Worksheets.Item("P600")
$ws.Range("E2").Value = 45096.48990774305
$ws.Range("E3").Value = 45096.48990774305
$ws.Range("E4").Value = 45096.48990774305
$ws.Range("E5").Value = 45096.48990774305

$ws = $wb.Worksheets.Item("IK91")
$ws.Range("E2").Value = 45096.48995410879
$ws.Range("E3").Value = 45096.48995410879
$ws.Range("E4").Value = 45096.48995410879

$ws = $wb.Worksheets.Item("IK92")
$ws.Range("E2").Value = 45096.49000885417
$ws.Range("E3").Value = 45096.49000885417
$ws.Range("E4").Value = 45096.49000898148

$ws = $wb.Worksheets.Item("IK93")
$ws.Range("E2").Value = 45096.49013689814
$ws.Range("E3").Value = 45096.49013703704
$ws.Range("E4").Value = 45096.49013703704
$ws.Range("E5").Value = 45096.49013689814
$ws.Range("E6").Value = 45096.49013689814

$ws = $wb.Worksheets.Item("IK94")
$ws.Range("E2").Value = 45096.4902321412
$ws.Range("E3").Value = 45096.49023223379
$ws.Range("E4").Value = 45096.49023223379
$ws.Range("E5").Value = 45096.4902321412

$ws = $wb.Worksheets.Item("IK95")
$ws.Range("E2").Value = 45096.49025758102

$ws = $wb.Worksheets.Item("IK96")
$ws.Range("E2").Value = 45096.49029898148
$ws.Range("E3").Value = 45096.49029907407
$ws.Range("E4").Value = 45096.49029898148

$ws = $wb.Worksheets.Item("IK101")
$ws.Range("E2").Value = 45096.49035333333
$ws.Range("E3").Value = 45096.4903534375
$ws.Range("E4").Value = 45096.4903534375

$ws = $wb.Worksheets.Item("IK102")
$ws.Range("E2").Value = 45096.4904037963
$ws.Range("E3").Value = 45096.4904037963
$ws.Range("E4").Value = 45096.4904037963
$ws.Range("E5").Value = 45096.4904037963

$ws = $wb.Worksheets.Item("IK103")
$ws.Range("E2").Value = 45096.49057321759
$ws.Range("E3").Value = 45096.49057332176
$ws.Range("E4").Value = 45096.49057332176
$ws.Range("E5").Value = 45096.49057332176
$ws.Range("E6").Value = 45096.49057332176

$ws = $wb.Worksheets.Item("IK104")
$ws.Range("E2").Value = 45096.49059289352
$ws.Range("E3").Value = 45096.49059289352
$ws.Range("E4").Value = 45096.49059289352
$ws.Range("E5").Value = 45096.49059289352

$ws = $wb.Worksheets.Item("IK105")
$ws.Range("E2").Value = 45096.49073003473
$ws.Range("E3").Value = 45096.49073011574
$ws.Range("E4").Value = 45096.49073011574
$ws.Range("E5").Value = 45096.49073003473

